$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 5952.5713
$ws.Range("I132").Value = 6225.95
$ws.Range("K132").Value = 18677.85
$ws.Range("M132").Value = -16147.85

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11754.825
$ws.Range("I32").Value = 6485.661
$ws.Range("J32").Value = 89475
$ws.Range("K32").Value = 6485.661
$ws.Range("L32").Value = 89475
$ws.Range("M32").Value = -6198.661
$ws.Range("N32").Value = -90049

# Row 37
$ws.Range("H37").Value = 9389.111000000001
$ws.Range("J37").Value = 9937.75
$ws.Range("L37").Value = 9937.75
$ws.Range("N37").Value = -10483.75

# Row 38
$ws.Range("H38").Value = 300
$ws.Range("I38").Value = 300
$ws.Range("K38").Value = 300
$ws.Range("M38").Value = 167

# Row 39
$ws.Range("H39").Value = 3500
$ws.Range("I39").Value = 2250
$ws.Range("J39").Value = 6000
$ws.Range("K39").Value = 2250
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = -1730
$ws.Range("N39").Value = -7040

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("L55").ClearContents()

# Row 61
$ws.Range("H61").Value = 3176211
$ws.Range("I61").Value = 3969804.2
$ws.Range("J61").Value = 1838.5714
$ws.Range("K61").Value = 3969804.2
$ws.Range("L61").Value = 1838.5714
$ws.Range("M61").Value = -3969592.2
$ws.Range("N61").Value = -2262.5714

# Row 63
$ws.Range("H63").Value = 4351.75
$ws.Range("I63").Value = 1598
$ws.Range("J63").Value = 6004
$ws.Range("K63").Value = 1598
$ws.Range("L63").Value = 6004
$ws.Range("M63").Value = -912
$ws.Range("N63").Value = -7376

# Row 66
$ws.Range("H66").Value = 4351.75
$ws.Range("I66").Value = 1598
$ws.Range("J66").Value = 6004
$ws.Range("K66").Value = 7990
$ws.Range("L66").Value = 30020
$ws.Range("M66").Value = -4558
$ws.Range("N66").Value = -36884

# Row 80
$ws.Range("H80").Value = 15200
$ws.Range("J80").Value = 19933.334
$ws.Range("L80").Value = 19933.334
$ws.Range("N80").Value = -21929.334

# Row 83
$ws.Range("H83").Value = 15200
$ws.Range("J83").Value = 19933.334
$ws.Range("L83").Value = 59800.00199999999
$ws.Range("N83").Value = -69784.00199999999

# Row 132
$ws.Range("H132").Value = 711687.7
$ws.Range("I132").Value = 987255.7
$ws.Range("J132").Value = 68695.664
$ws.Range("K132").Value = 2961767.1
$ws.Range("L132").Value = 206086.992
$ws.Range("M132").Value = -2959237.1
$ws.Range("N132").Value = -211146.992

# Row 136
$ws.Range("H136").Value = 3176211
$ws.Range("I136").Value = 3969804.2
$ws.Range("J136").Value = 1838.5714
$ws.Range("K136").Value = 11909412.6
$ws.Range("L136").Value = 5515.7142
$ws.Range("M136").Value = -11906862.6
$ws.Range("N136").Value = -10615.7142

$ws = $wb.Worksheets.Item("BSM")
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("L45").ClearContents()

# Row 82
$ws.Range("H82").Value = 17008.143
$ws.Range("I82").Value = 10314.25
$ws.Range("K82").Value = 10314.25
$ws.Range("M82").Value = -9931.25

# Row 85
$ws.Range("H85").Value = 17008.143
$ws.Range("I85").Value = 10314.25
$ws.Range("K85").Value = 10314.25
$ws.Range("M85").Value = -8988.25

# Row 134
$ws.Range("H134").Value = 11168006
$ws.Range("I134").Value = 12563638
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 37690914
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -37688379
$ws.Range("N134").Value = -13920

$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 10546.5
$ws.Range("I50").Value = 9181.5
$ws.Range("J50").Value = 10887.75
$ws.Range("K50").Value = 9181.5
$ws.Range("L50").Value = 10887.75
$ws.Range("M50").Value = -8556.5
$ws.Range("N50").Value = -12137.75

# Row 51
$ws.Range("H51").Value = 9975.375
$ws.Range("I51").Value = 9050
$ws.Range("J51").Value = 10283.833
$ws.Range("K51").Value = 9050
$ws.Range("L51").Value = 10283.833
$ws.Range("M51").Value = -8314
$ws.Range("N51").Value = -11755.833

# Row 58
$ws.Range("H58").Value = 2916.5898
$ws.Range("I58").Value = 1691.2354
$ws.Range("J58").Value = 3863.4546
$ws.Range("K58").Value = 1691.2354
$ws.Range("L58").Value = 3863.4546
$ws.Range("M58").Value = -1488.2354
$ws.Range("N58").Value = -4269.4546

# Row 59
$ws.Range("H59").Value = 16809.8
$ws.Range("J59").Value = 16809.8
$ws.Range("L59").Value = 16809.8
$ws.Range("N59").Value = -19099.8

# Row 60
$ws.Range("H60").Value = 9308.75
$ws.Range("J60").Value = 9970.5
$ws.Range("L60").Value = 9970.5
$ws.Range("N60").Value = -10992.5

# Row 61
$ws.Range("H61").Value = 9975.375
$ws.Range("I61").Value = 9050
$ws.Range("J61").Value = 10283.833
$ws.Range("K61").Value = 9050
$ws.Range("L61").Value = 10283.833
$ws.Range("M61").Value = -8702
$ws.Range("N61").Value = -10979.833

# Row 74
$ws.Range("H74").Value = 14622.667
$ws.Range("J74").Value = 16200.5
$ws.Range("L74").Value = 16200.5
$ws.Range("N74").Value = -17948.5

# Row 77
$ws.Range("H77").Value = 14622.667
$ws.Range("J77").Value = 16200.5
$ws.Range("L77").Value = 48601.5
$ws.Range("N77").Value = -57337.5

# Row 99
$ws.Range("H99").Value = 1484
$ws.Range("I99").Value = 1349.9
$ws.Range("J99").Value = 1651.625
$ws.Range("K99").Value = 1349.9
$ws.Range("L99").Value = 1651.625
$ws.Range("M99").Value = 148.0999999999999
$ws.Range("N99").Value = -4647.625

# Row 126
$ws.Range("H126").Value = 1484
$ws.Range("I126").Value = 1349.9
$ws.Range("J126").Value = 1651.625
$ws.Range("K126").Value = 4049.7
$ws.Range("L126").Value = 4954.875
$ws.Range("M126").Value = -1579.7
$ws.Range("N126").Value = -9894.875

# Row 132
$ws.Range("H132").Value = 4476.9165
$ws.Range("I132").Value = 5419.2856
$ws.Range("J132").Value = 3157.6
$ws.Range("K132").Value = 16257.8568
$ws.Range("L132").Value = 9472.799999999999
$ws.Range("M132").Value = -13727.8568
$ws.Range("N132").Value = -14532.8

# Row 134
$ws.Range("H134").Value = 3037.375
$ws.Range("I134").Value = 2965.8333
$ws.Range("J134").Value = 3252
$ws.Range("K134").Value = 8897.499899999999
$ws.Range("L134").Value = 9756
$ws.Range("M134").Value = -6362.499899999999
$ws.Range("N134").Value = -14826

# Row 136
$ws.Range("H136").Value = 2916.5898
$ws.Range("I136").Value = 1691.2354
$ws.Range("J136").Value = 3863.4546
$ws.Range("K136").Value = 5073.706200000001
$ws.Range("L136").Value = 11590.3638
$ws.Range("M136").Value = -2523.706200000001
$ws.Range("N136").Value = -16690.3638

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 13158845
$ws.Range("I5").Value = 867.3333
$ws.Range("J5").Value = 15625966
$ws.Range("K5").Value = 2601.9999
$ws.Range("L5").Value = 46877898
$ws.Range("M5").Value = -2489.9999
$ws.Range("N5").Value = -46878122

# Row 131
$ws.Range("H131").Value = 2163.28
$ws.Range("I131").Value = 50000
$ws.Range("J131").Value = 1516.8379
$ws.Range("K131").Value = 150000
$ws.Range("L131").Value = 4550.5137
$ws.Range("M131").Value = -144960
$ws.Range("N131").Value = -14630.5137

# Row 135
$ws.Range("H135").Value = 13158845
$ws.Range("I135").Value = 867.3333
$ws.Range("J135").Value = 15625966
$ws.Range("K135").Value = 7805.9997
$ws.Range("L135").Value = 140633694
$ws.Range("M135").Value = -5270.9997
$ws.Range("N135").Value = -140638764

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 11022.489
$ws.Range("I122").Value = 2175.6843
$ws.Range("J122").Value = 17487.46
$ws.Range("K122").Value = 6527.0529
$ws.Range("L122").Value = 52462.38
$ws.Range("M122").Value = -4077.0529
$ws.Range("N122").Value = -57362.38

# Row 132
$ws.Range("H132").Value = 2028.8823
$ws.Range("I132").Value = 1834.9286
$ws.Range("J132").Value = 2164.65
$ws.Range("K132").Value = 5504.7858
$ws.Range("L132").Value = 6493.950000000001
$ws.Range("M132").Value = -2974.7858
$ws.Range("N132").Value = -11553.95

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1429.8334
$ws.Range("I40").Value = 1419.8235
$ws.Range("K40").Value = 1419.8235
$ws.Range("M40").Value = -1283.8235

# Row 136
$ws.Range("H136").Value = 2091.3635
$ws.Range("I136").Value = 1595.2632
$ws.Range("J136").Value = 5233.3335
$ws.Range("K136").Value = 4785.7896
$ws.Range("L136").Value = 15700.0005
$ws.Range("M136").Value = -2235.7896
$ws.Range("N136").Value = -20800.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1470.5333
$ws.Range("I122").Value = 1470.5333
$ws.Range("K122").Value = 4411.5999
$ws.Range("M122").Value = -1961.5999
